# Applies the row-data corrections for rows 14, 16, 17, 18 and 19
# (swapped/corrected species records) as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 14 ----
$ws.Range("A14").Value = 111798755
$ws.Range("B14").Value = 90709
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 5448
$ws.Range("F14").Value = "Svartvit taggsvamp"
$ws.Range("G14").Value = "Phellodon connatus"
$ws.Range("H14").Value = "(Schultz) nom.prov"
$ws.Range("Q14").Value = 753031
$ws.Range("R14").Value = 7090921
$ws.Range("S14").Value = 25
$ws.Range("AF14").ClearContents()
$ws.Range("AI14").ClearContents()

# ---- Row 16 ----
$ws.Range("A16").Value = 111798757
$ws.Range("B16").Value = 81076
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 5046
$ws.Range("F16").Value = "Grön jordtunga"
$ws.Range("G16").Value = "Microglossum viride"
$ws.Range("H16").Value = "(Pers.:Fr.) Gillet"
$ws.Range("Q16").Value = 753109
$ws.Range("R16").Value = 7091008
$ws.Range("S16").Value = 100
$ws.Range("AF16").Value = "mikroskoperad"
$ws.Range("AI16").Value = "Granskog"

# ---- Row 17 ----
$ws.Range("A17").Value = 111961472
$ws.Range("B17").Value = 90843
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 5448
$ws.Range("F17").Value = "Svartvit taggsvamp"
$ws.Range("G17").Value = "Phellodon connatus"
$ws.Range("H17").Value = "(Schultz) nom.prov"
$ws.Range("I17").ClearContents()
$ws.Range("J17").ClearContents()
$ws.Range("AF17").ClearContents()
$ws.Range("AX17").Value = "Stefan Phalagorn Bergström, Annika  Carlberg , Andreas Estensen, Ola Elleström, Anne Järvinen, Emma Sewell, Thomas Strid"

# ---- Row 18 ----
$ws.Range("A18").Value = 111961716
$ws.Range("B18").Value = 81193
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 5046
$ws.Range("F18").Value = "Grön jordtunga"
$ws.Range("G18").Value = "Microglossum viride"
$ws.Range("H18").Value = "(Pers.:Fr.) Gillet"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "2"
$ws.Range("J18").Value = "mycel"
$ws.Range("AF18").Value = "mikroskoperad"
$ws.Range("AX18").Value = "Stefan Phalagorn Bergström, Andreas Estensen, Annika  Carlberg , Ola Elleström, Thomas Strid, Anne Järvinen, Emma Sewell"

# ---- Row 19 ----
$ws.Range("B19").Value = 89936
